$wb = $excel.ActiveWorkbook

$wsPresupuesto = $wb.Worksheets.Item("Presupuesto")
$wsCompras = $wb.Worksheets.Item("Compras")
$wsInventario = $wb.Worksheets.Item("Inventario_cero_coste")

# --- Compras: new purchase row (row 4) -----------------------------------
# Add the hyperlink first (links B4 back to the inventory item), then set
# its value -- Hyperlinks.Add(..) would otherwise stomp the cell's value
# with the link's display text.
$wsCompras.Hyperlinks.Add($wsCompras.Range("B4"), "", "Inventario_cero_coste!B11", "", "Inventario_cero_coste!B11")
$wsCompras.Range("B4").Value = 9
$wsCompras.Range("C4").Value = "Manu"
$wsCompras.Range("D4").Value = 41.73

# --- Inventario_cero_coste: new stepper-motor inventory row (row 11) -----
# Add the new item's amazon hyperlink before filling in the other cells so
# the new shared strings land in the same order as the authored workbook.
$wsInventario.Hyperlinks.Add($wsInventario.Range("I11"), "https://www.amazon.es/dp/B07SWYFCQV?smid=A3LC78H97WEBWA&ref_=chk_typ_imgToDp&th=1")

$wsInventario.Range("C11").Value = "Motor 17HS4401"
$wsInventario.Range("E11").Value = "Motor, Stepper"
$wsInventario.Range("F11").Value = "Motor Paso Paso Nema 17"
$wsInventario.Range("G11").Value = "Manu"

# Existing rows 8 & 9 already showed their source URLs as plain text --
# turn those into real clickable hyperlinks too.
$wsInventario.Hyperlinks.Add($wsInventario.Range("I8"), "https://servodatabase.com/servo/springrc/sm-s2309s")
$wsInventario.Hyperlinks.Add($wsInventario.Range("I9"), "https://www.digipart.com/part/3590S-2-103?utm_source=bing&utm_medium=cpc&utm_campaign=arrow&utm_term=3590S-2-103&utm_content=ad_arrow_b")

# --- Selections on each sheet, restoring Presupuesto as the active tab ---
$wsCompras.Range("D5").Select()
$wsInventario.Range("A3").Select()
$wsPresupuesto.Activate()
